# Apply "Improved container source feature names and descriptions" edit.
$wb = $excel.ActiveWorkbook

$wsContainer = $wb.Worksheets.Item("Container Features")
$wsComponent = $wb.Worksheets.Item("Component Features")
$wsExample   = $wb.Worksheets.Item("Example Features")

# --- Container Features sheet edits ---
# Row 2 (Sphinx Runtime) description: "sources" -> "source code"
$wsContainer.Range("C2").Value = "Includes Sphinx runtime binaries (without source code and documentation)"

# Row 5 (Sphinx SDK) description: "sources" -> "source code"
$wsContainer.Range("C5").Value = "Includes Sphinx runtime, source code and documentation (but no examples)"

# Row 6 (Sphinx Examples) description: "(with sources)" -> "(with source code)"
$wsContainer.Range("C6").Value = "Includes Sphinx examples (with source code)"

# Row 3 (was "Sphinx Sources") renamed feature + description to "automatically generated"
$wsContainer.Range("B3").Value = "automatically generated"
$wsContainer.Range("C3").Value = "automatically generated"

# Row 3 no longer needs its taller custom height now that the text shrank -
# let Excel recompute the (now default) auto row height.
$wsContainer.Rows.Item(3).AutoFit()

# Update selections to match final authored view state
$wsContainer.Range("C6").Select()
$wsComponent.Range("C13").Select()
$wsContainer.Activate()
